$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections / new content ---
# These edits are applied in the same order the original author made them so
# that the shared-strings table ends up laid out identically.

# output_dir description: added trailing period
$ws.Range("BM1").Value = "The output directory for this project. All the results from the pipeline as well as the intermediary steps will be placed here."

# project_long_name description: "sample project" -> "same project"
$ws.Range("F1").Value = "You are free to describe the project name with any characters you want. Within the same project however this field must be identical for all entries."

# --- New columns: phred_window_size (BW) and phred_threshold (BX) ---
$ws.Range("BW2").Value = "phred_window_size"
$ws.Range("BX2").Value = "phred_threshold"

# max_read_len description: "tunrcation" -> "truncation"
$ws.Range("BV1").Value = "How long can a read be at a maximum (after primer truncation) before we discard that read? No value will disable this check."

# min_read_len description: "tunrcation" -> "truncation"
$ws.Range("BU1").Value = "How short can a read be at a minimum (after primer truncation) before we discard that read? No value will disable this check."

$ws.Range("BW1").Value = "How large should the rolling average window be in base pairs when we filter reads based on their quality scores? No value will disable this check."
$ws.Range("BX1").Value = "What quality score value should every window have at a minimum before we discard that read? No value will disable this check."

# primer_max_dist description: appended "No value will disable this check."
$ws.Range("BT1").Value = "How far away from the start of the read can we locate the primer sequence before we discard that read? Enter a value between 0 and 150. No value will disable this check."

# Optional-metadata free-text description: removed duplicated space
$ws.Range("BK1").Value = "You can add any extra columns of optional metadata here."

# used/"If you want to exclude..." description: added trailing period
$ws.Range("B1").Value = "If you want to exclude a particular sample from analysis enter ""no"", otherwise enter ""yes""."

# New column widths to match the newly-inserted columns
$ws.Columns("BW").ColumnWidth = 24.333333333333332
$ws.Columns("BX").ColumnWidth = 21

# Restore the scroll/selection state of the sheet view
$ws.Range("F13").Select()
